# Fix a minor mistake:
#  1. Append a separate "1" run after "Story ID 1" (-> "Story ID 1" + "1") and
#     move the "_GoBack" bookmark to sit right after it.
#  2. Merge the "NOT " / "IMPLEMENTED " runs and the " " / "INSURANCE..." runs
#     into single runs each.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part A: add a new run containing "1" right after the "Story ID 1" heading,
# and relocate the "_GoBack" bookmark there.
# ---------------------------------------------------------------------------

# Build the extra "1" text on a scratch range at the very end of the story
# (so FormattedText carries no special formatting), then paste it in via
# FormattedText so it lands as its own <w:r> instead of being merged into the
# neighbouring run the way InsertAfter/TypeText would.
$endOfDoc = $d.Content.Text.Length - 1
$scratch = $d.Range($endOfDoc, $endOfDoc)
$scratch.InsertAfter("1")
$scratchRange = $d.Range($endOfDoc, $endOfDoc + 1)
$ft = $scratchRange.FormattedText

$headingPara = $d.Paragraphs(2)
$headingEnd = $headingPara.Range.End
$insertPoint = $d.Range($headingEnd - 1, $headingEnd - 1)
$insertPoint.FormattedText = $ft

# Remove the scratch copy we typed at the end of the document.
$tailLen = $d.Content.Text.Length
$scratchRange2 = $d.Range($tailLen - 2, $tailLen - 1)
$scratchRange2.Delete()

# Work out where the new run ends (right after the "1" we just added).
$headingPara = $d.Paragraphs(2)
$bookmarkPos = $headingPara.Range.End - 1

# The engine mis-places a bookmark added exactly at "paragraph end minus one"
# (it always snaps to (0,28) instead of the requested collapsed position), so
# add a temporary trailing character to move that boundary out of the way,
# add/relocate the bookmark next to the real text, then remove the helper
# character again.
$tempAnchor = $d.Range($bookmarkPos, $bookmarkPos)
$tempAnchor.InsertAfter("Z")

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$helperRange = $d.Range($bookmarkPos, $bookmarkPos + 1)
$helperRange.Delete()

# ---------------------------------------------------------------------------
# Part B: coalesce the split runs in the "NOT IMPLEMENTED ... INSURANCE ..."
# paragraph. Running Find/Replace over the text re-normalises adjacent runs
# that share identical formatting into a single run.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("NOT IMPLEMENTED ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "NOT IMPLEMENTED ", 2)
